$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SKU codes (PRTS -> QRTS) and product names for rows 2-6
$ws.Range("A2").Value = "QRTS00001"
$ws.Range("A3").Value = "QRTS00002"
$ws.Range("A4").Value = "QRTS00003"
$ws.Range("B4").Value = "Alcor Deep BEEF Fry Wok"
$ws.Range("A5").Value = "QRTS00004"
$ws.Range("B5").Value = "Alcor Shallow KOBE Work"
$ws.Range("A6").Value = "QRTS00005"
$ws.Range("B6").Value = "Alcor Deep MICIN Fry Wok"

# Move the active selection to B2
$ws.Range("B2").Select()

# Remove the duplicate-values conditional formatting rule on column A
$ws.Range("A:A").FormatConditions.Delete()
